$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of updates: row, new Price (D), new Volume(1h) (E)
# Rows whose Price (D) is unchanged use $null as a "skip" marker.
$updates = @(
    @{ Row = 2;  D = "27.927.05";   E = "  +0.82%  " },
    @{ Row = 3;  D = "1.880.08";    E = "  +0.35%  " },
    @{ Row = 4;  D = "1.018";       E = "  +1.37%  " },
    @{ Row = 5;  D = "334.97";      E = "  +0.87%  " },
    @{ Row = 6;  D = $null;         E = "  +1.22%  " },
    @{ Row = 7;  D = "0.4684";      E = "  -0.75%  " },
    @{ Row = 8;  D = "0.3913";      E = "  -0.99%  " },
    @{ Row = 9;  D = "46.87";       E = "  -1.30%  " },
    @{ Row = 10; D = "0.07944";     E = "  -1.09%  " },
    @{ Row = 11; D = "1.007";       E = "  -1.46%  " },
    @{ Row = 12; D = "21.57";       E = "  -1.06%  " },
    @{ Row = 13; D = "1.889.43";    E = "  +0.42%  " },
    @{ Row = 14; D = "5.949";       E = "  -0.10%  " },
    @{ Row = 15; D = "7.091";       E = "  -0.73%  " },
    @{ Row = 16; D = "1.020";       E = "  +1.24%  " },
    @{ Row = 17; D = "0.06785";     E = "  +2.33%  " },
    @{ Row = 18; D = "87.50";       E = "  +0.35%  " },
    @{ Row = 19; D = "0.00001044";  E = "  -0.01%  " },
    @{ Row = 20; D = "17.01";       E = "  -0.95%  " },
    @{ Row = 22; D = "27.941.85";   E = "  +0.81%  " },
    @{ Row = 23; D = "5.467";       E = "  -0.50%  " },
    @{ Row = 24; D = "10.96";       E = "  -0.45%  " },
    @{ Row = 25; D = "2.359";       E = "  +2.55%  " },
    @{ Row = 26; D = "2.111.12";    E = "  +0.28%  " },
    @{ Row = 27; D = "159.44";      E = "  +1.98%  " },
    @{ Row = 28; D = "19.89";       E = "  -1.48%  " },
    @{ Row = 29; D = "2.069";       E = "  -1.35%  " },
    @{ Row = 30; D = "5.454";       E = "  -1.96%  " },
    @{ Row = 31; D = "120.80";      E = "  -1.33%  " },
    @{ Row = 32; D = "0.09526";     E = "  -0.45%  " },
    @{ Row = 33; D = "0.9551";      E = "  -1.25%  " },
    @{ Row = 34; D = "3.661";       E = "  +0.79%  " },
    @{ Row = 35; D = "5.318";       E = "  +0.38%  " },
    @{ Row = 36; D = "1.348";       E = "  -7.15%  " },
    @{ Row = 37; D = "0.06113";     E = "  -0.04%  " },
    @{ Row = 38; D = "0.02234";     E = "  -1.12%  " },
    @{ Row = 39; D = "1.201";       E = "  -1.87%  " },
    @{ Row = 40; D = $null;         E = "  +1.23%  " },
    @{ Row = 41; D = "8.091";       E = "  -0.92%  " },
    @{ Row = 42; D = "0.5880";      E = "  -1.70%  " },
    @{ Row = 43; D = "0.1893";      E = "  -0.58%  " },
    @{ Row = 44; D = "10.16";       E = "  -1.02%  " },
    @{ Row = 45; D = "1.274";       E = "  +1.86%  " },
    @{ Row = 46; D = "0.5636";      E = "  -0.87%  " },
    @{ Row = 47; D = $null;         E = "  -0.89%  " },
    @{ Row = 48; D = "3.398";       E = "  -0.23%  " },
    @{ Row = 49; D = "1.915";       E = "  -0.73%  " },
    @{ Row = 50; D = "0.06857";     E = "  +0.56%  " },
    @{ Row = 51; D = "113.54";      E = "  +0.94%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        # Force text interpretation so numeric-looking price strings (e.g.
        # "1.020") keep their exact digits instead of being coerced into a
        # number (which would drop significant trailing zeros).
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
